$wb = $excel.ActiveWorkbook

# --- Service Contacts sheet: selection change + new column width ---
$wsServiceContacts = $wb.Worksheets.Item("Service Contacts")
$wsServiceContacts.Range("D3").Select() | Out-Null
$wsServiceContacts.Columns.Item(1).ColumnWidth = 13.666666666666666

# --- Practitioners sheet: new column widths, new data row, selection change ---
$wsPractitioners = $wb.Worksheets.Item("Practitioners")

$wsPractitioners.Columns.Item(1).ColumnWidth = 13.833333333333334
$wsPractitioners.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsPractitioners.Columns.Item(6).ColumnWidth = 12.0

$wsPractitioners.Range("A6").Value = "PHN999:NFP02"
$wsPractitioners.Range("B6").Value = "P01"
$wsPractitioners.Range("C6").Value = 8
$wsPractitioners.Range("D6").Value = 1
$wsPractitioners.Range("E6").Value = 1973
$wsPractitioners.Range("F6").Value = 2
$wsPractitioners.Range("G6").Value = 1
$wsPractitioners.Range("H6").Value = 1
$wsPractitioners.Range("I6").Value = "tag1"

$wsPractitioners.Range("G1:G1048576").Select() | Out-Null

# Restore original active sheet ("Episodes", activeTab=3) so workbook-level
# active tab / tabSelected bookkeeping is unaffected by the edits above.
$wb.Worksheets.Item("Episodes").Activate() | Out-Null
